$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the data row for matricula GO492007807 (TIENE GOMES DE LIMA ABREU),
# which lives at row 724 of the "Export" sheet.
$ws.Rows.Item(724).Delete()
